$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values for columns A and B, rows 5-29 (row, colA, colB)
$data = @(
    @(5,  2759.7575790000001, 249.23779400000001),
    @(6,  2758.560023,        248.0304691),
    @(7,  2758.4131649999999, 247.8835297),
    @(8,  2759.8135520000001, 249.27926350000001),
    @(9,  2759.5015239999998, 248.9732731),
    @(10, 2759.5711019999999, 249.0488517),
    @(11, 2759.0370320000002, 248.50312690000001),
    @(12, 2759.7899600000001, 249.2539573),
    @(13, 2758.291995,        247.7483781),
    @(14, 2759.3015230000001, 248.77817289999999),
    @(15, 2758.8650579999999, 248.3392187),
    @(16, 2759.2593459999998, 248.7297374),
    @(17, 2758.38805,         247.85348870000001),
    @(18, 2759.5090070000001, 248.9718881),
    @(19, 2760.1874859999998, 249.66849780000001),
    @(20, 2758.3600219999998, 247.8151484),
    @(21, 2759.829416,        249.30228399999999),
    @(22, 2757.5006990000002, 246.9629544),
    @(23, 2758.6724049999998, 248.14187089999999),
    @(24, 2759.6117279999999, 249.0811946),
    @(25, 2759.6416599999998, 249.1082696),
    @(26, 2757.511039,        246.9621108),
    @(27, 2759.2797540000001, 248.74010480000001),
    @(28, 2754.879735,        244.50743370000001),
    @(29, 2759.8348590000001, 249.304597)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $a = $entry[1]
    $b = $entry[2]
    $ws.Cells.Item($r, 1).Value = $a
    $ws.Cells.Item($r, 2).Value = $b
}

# Move the active cell selection to B4, as seen in the edited workbook.
$ws.Range("B4").Select()
